$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "FAILED: $old"
    }
}

# 1. Capitalize first word of the document
ReplaceText "put together a series" "Put together a series"

# 2. Remove "book"/"website" hyperlinks and surrounding parenthetical
ReplaceText "that influence. And this is where mediation analysis can be useful. (If you want to delve deeply into the topic, I recommend you check out this book by Tyler VanderWeele, or this nice website developed at Columbia University.)" "that influence. And this is where mediation analysis can be useful. "

# 4. Remove "paper" hyperlink and surrounding parenthetical
ReplaceText "(or quantities) that arise in a mediation analysis. (I draw on a paper by Imai, Keele and Yamamoto for the terminology, as there is not complete agreement on what to call various quantities. The estimation methods and software used here are also described in the paper.)" "(or quantities) that arise in a mediation analysis. "

# 5. Remove "In an earlier post, I described the concept of potential outcomes. "
ReplaceText "In an earlier post, I described the concept of potential outcomes. I extend" "I extend"

# 8. Remove "(see Imai et al for the details), "
ReplaceText "independence, and consistency (see Imai et al for the details), the average causal" "independence, and consistency, the average causal"

# 16. Remove "(again, Imai et al provide these), "
ReplaceText "I will not go into the important details here (again, Imai et al provide these), but here are formulas" "I will not go into the important details here, but here are formulas"
